$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple isolated edits (C/D/E only) ---
$ws.Cells.Item(5095, 3).Value = 2431
$ws.Cells.Item(5095, 4).Value = 39
$ws.Cells.Item(5095, 5).Value = 2470
$ws.Cells.Item(5172, 3).Value = 2772
$ws.Cells.Item(5172, 4).Value = 124
$ws.Cells.Item(5172, 5).Value = 2896
$ws.Cells.Item(5251, 3).Value = 2278
$ws.Cells.Item(5251, 4).Value = 103
$ws.Cells.Item(5251, 5).Value = 2381
$ws.Cells.Item(5329, 3).Value = 1673
$ws.Cells.Item(5329, 4).Value = 68
$ws.Cells.Item(5329, 5).Value = 1741
$ws.Cells.Item(5357, 3).Value = 1079
$ws.Cells.Item(5357, 4).Value = 27
$ws.Cells.Item(5357, 5).Value = 1106
$ws.Cells.Item(5407, 3).Value = 1676
$ws.Cells.Item(5407, 4).Value = 61
$ws.Cells.Item(5407, 5).Value = 1737
$ws.Cells.Item(5425, 3).Value = 788
$ws.Cells.Item(5425, 4).Value = 36
$ws.Cells.Item(5425, 5).Value = 824
$ws.Cells.Item(5591, 3).Value = 1157
$ws.Cells.Item(5591, 4).Value = 50
$ws.Cells.Item(5591, 5).Value = 1207
$ws.Cells.Item(5603, 3).Value = 468
$ws.Cells.Item(5603, 4).Value = 55
$ws.Cells.Item(5603, 5).Value = 523
$ws.Cells.Item(5638, 3).Value = 73
$ws.Cells.Item(5638, 4).Value = 3
$ws.Cells.Item(5638, 5).Value = 76
$ws.Cells.Item(5647, 3).Value = 1037
$ws.Cells.Item(5647, 4).Value = 33
$ws.Cells.Item(5647, 5).Value = 1070
$ws.Cells.Item(5654, 3).Value = 1323
$ws.Cells.Item(5654, 4).Value = 38
$ws.Cells.Item(5654, 5).Value = 1361

# --- Insert 9 new rows at the end of the 2021-01-19 block (after row 5723) ---
$ws.Rows("5724:5732").Insert()

# --- Rewrite the full 2021-01-19 block (rows 5665-5732) ---
$ws.Cells.Item(5665, 1).Value = 44215
$ws.Cells.Item(5665, 2).Value = 'Okres Bánovce nad Bebravou'
$ws.Cells.Item(5665, 3).Value = 0
$ws.Cells.Item(5665, 4).Value = 1
$ws.Cells.Item(5665, 5).Value = 1
$ws.Cells.Item(5666, 1).Value = 44215
$ws.Cells.Item(5666, 2).Value = 'Okres Banská Bystrica'
$ws.Cells.Item(5666, 3).Value = 1823
$ws.Cells.Item(5666, 4).Value = 37
$ws.Cells.Item(5666, 5).Value = 1860
$ws.Cells.Item(5667, 1).Value = 44215
$ws.Cells.Item(5667, 2).Value = 'Okres Banská Štiavnica'
$ws.Cells.Item(5667, 3).Value = 378
$ws.Cells.Item(5667, 4).Value = 23
$ws.Cells.Item(5667, 5).Value = 401
$ws.Cells.Item(5668, 1).Value = 44215
$ws.Cells.Item(5668, 2).Value = 'Okres Bardejov'
$ws.Cells.Item(5668, 3).Value = 86
$ws.Cells.Item(5668, 4).Value = 1
$ws.Cells.Item(5668, 5).Value = 87
$ws.Cells.Item(5669, 1).Value = 44215
$ws.Cells.Item(5669, 2).Value = 'Okres Bratislava I'
$ws.Cells.Item(5669, 3).Value = 232
$ws.Cells.Item(5669, 4).Value = 5
$ws.Cells.Item(5669, 5).Value = 237
$ws.Cells.Item(5670, 1).Value = 44215
$ws.Cells.Item(5670, 2).Value = 'Okres Bratislava II'
$ws.Cells.Item(5670, 3).Value = 1542
$ws.Cells.Item(5670, 4).Value = 78
$ws.Cells.Item(5670, 5).Value = 1620
$ws.Cells.Item(5671, 1).Value = 44215
$ws.Cells.Item(5671, 2).Value = 'Okres Bratislava III'
$ws.Cells.Item(5671, 3).Value = 879
$ws.Cells.Item(5671, 4).Value = 22
$ws.Cells.Item(5671, 5).Value = 901
$ws.Cells.Item(5672, 1).Value = 44215
$ws.Cells.Item(5672, 2).Value = 'Okres Bratislava IV'
$ws.Cells.Item(5672, 3).Value = 252
$ws.Cells.Item(5672, 4).Value = 17
$ws.Cells.Item(5672, 5).Value = 269
$ws.Cells.Item(5673, 1).Value = 44215
$ws.Cells.Item(5673, 2).Value = 'Okres Bratislava V'
$ws.Cells.Item(5673, 3).Value = 365
$ws.Cells.Item(5673, 4).Value = 16
$ws.Cells.Item(5673, 5).Value = 381
$ws.Cells.Item(5674, 1).Value = 44215
$ws.Cells.Item(5674, 2).Value = 'Okres Brezno'
$ws.Cells.Item(5674, 3).Value = 29
$ws.Cells.Item(5674, 4).Value = 5
$ws.Cells.Item(5674, 5).Value = 34
$ws.Cells.Item(5675, 1).Value = 44215
$ws.Cells.Item(5675, 2).Value = 'Okres Bytča'
$ws.Cells.Item(5675, 3).Value = 303
$ws.Cells.Item(5675, 4).Value = 31
$ws.Cells.Item(5675, 5).Value = 334
$ws.Cells.Item(5676, 1).Value = 44215
$ws.Cells.Item(5676, 2).Value = 'Okres Čadca'
$ws.Cells.Item(5676, 3).Value = 1094
$ws.Cells.Item(5676, 4).Value = 19
$ws.Cells.Item(5676, 5).Value = 1113
$ws.Cells.Item(5677, 1).Value = 44215
$ws.Cells.Item(5677, 2).Value = 'Okres Detva'
$ws.Cells.Item(5677, 3).Value = 31
$ws.Cells.Item(5677, 4).Value = 0
$ws.Cells.Item(5677, 5).Value = 31
$ws.Cells.Item(5678, 1).Value = 44215
$ws.Cells.Item(5678, 2).Value = 'Okres Dolný Kubín'
$ws.Cells.Item(5678, 3).Value = 225
$ws.Cells.Item(5678, 4).Value = 10
$ws.Cells.Item(5678, 5).Value = 235
$ws.Cells.Item(5679, 1).Value = 44215
$ws.Cells.Item(5679, 2).Value = 'Okres Dunajská Streda'
$ws.Cells.Item(5679, 3).Value = 418
$ws.Cells.Item(5679, 4).Value = 59
$ws.Cells.Item(5679, 5).Value = 477
$ws.Cells.Item(5680, 1).Value = 44215
$ws.Cells.Item(5680, 2).Value = 'Okres Galanta'
$ws.Cells.Item(5680, 3).Value = 9
$ws.Cells.Item(5680, 4).Value = 0
$ws.Cells.Item(5680, 5).Value = 9
$ws.Cells.Item(5681, 1).Value = 44215
$ws.Cells.Item(5681, 2).Value = 'Okres Gelnica'
$ws.Cells.Item(5681, 3).Value = 1
$ws.Cells.Item(5681, 4).Value = 0
$ws.Cells.Item(5681, 5).Value = 1
$ws.Cells.Item(5682, 1).Value = 44215
$ws.Cells.Item(5682, 2).Value = 'Okres Hlohovec'
$ws.Cells.Item(5682, 3).Value = 12
$ws.Cells.Item(5682, 4).Value = 0
$ws.Cells.Item(5682, 5).Value = 12
$ws.Cells.Item(5683, 1).Value = 44215
$ws.Cells.Item(5683, 2).Value = 'Okres Humenné'
$ws.Cells.Item(5683, 3).Value = 261
$ws.Cells.Item(5683, 4).Value = 14
$ws.Cells.Item(5683, 5).Value = 275
$ws.Cells.Item(5684, 1).Value = 44215
$ws.Cells.Item(5684, 2).Value = 'Okres Ilava'
$ws.Cells.Item(5684, 3).Value = 280
$ws.Cells.Item(5684, 4).Value = 14
$ws.Cells.Item(5684, 5).Value = 294
$ws.Cells.Item(5685, 1).Value = 44215
$ws.Cells.Item(5685, 2).Value = 'Okres Kežmarok'
$ws.Cells.Item(5685, 3).Value = 179
$ws.Cells.Item(5685, 4).Value = 9
$ws.Cells.Item(5685, 5).Value = 188
$ws.Cells.Item(5686, 1).Value = 44215
$ws.Cells.Item(5686, 2).Value = 'Okres Komárno'
$ws.Cells.Item(5686, 3).Value = 489
$ws.Cells.Item(5686, 4).Value = 45
$ws.Cells.Item(5686, 5).Value = 534
$ws.Cells.Item(5687, 1).Value = 44215
$ws.Cells.Item(5687, 2).Value = 'Okres Košice I'
$ws.Cells.Item(5687, 3).Value = 54
$ws.Cells.Item(5687, 4).Value = 4
$ws.Cells.Item(5687, 5).Value = 58
$ws.Cells.Item(5688, 1).Value = 44215
$ws.Cells.Item(5688, 2).Value = 'Okres Košice II'
$ws.Cells.Item(5688, 3).Value = 808
$ws.Cells.Item(5688, 4).Value = 6
$ws.Cells.Item(5688, 5).Value = 814
$ws.Cells.Item(5689, 1).Value = 44215
$ws.Cells.Item(5689, 2).Value = 'Okres Košice IV'
$ws.Cells.Item(5689, 3).Value = 286
$ws.Cells.Item(5689, 4).Value = 6
$ws.Cells.Item(5689, 5).Value = 292
$ws.Cells.Item(5690, 1).Value = 44215
$ws.Cells.Item(5690, 2).Value = 'Okres Kysucké Nové Mesto'
$ws.Cells.Item(5690, 3).Value = 287
$ws.Cells.Item(5690, 4).Value = 14
$ws.Cells.Item(5690, 5).Value = 301
$ws.Cells.Item(5691, 1).Value = 44215
$ws.Cells.Item(5691, 2).Value = 'Okres Levice'
$ws.Cells.Item(5691, 3).Value = 1644
$ws.Cells.Item(5691, 4).Value = 156
$ws.Cells.Item(5691, 5).Value = 1800
$ws.Cells.Item(5692, 1).Value = 44215
$ws.Cells.Item(5692, 2).Value = 'Okres Levoča'
$ws.Cells.Item(5692, 3).Value = 213
$ws.Cells.Item(5692, 4).Value = 5
$ws.Cells.Item(5692, 5).Value = 218
$ws.Cells.Item(5693, 1).Value = 44215
$ws.Cells.Item(5693, 2).Value = 'Okres Liptovský Mikuláš'
$ws.Cells.Item(5693, 3).Value = 691
$ws.Cells.Item(5693, 4).Value = 20
$ws.Cells.Item(5693, 5).Value = 711
$ws.Cells.Item(5694, 1).Value = 44215
$ws.Cells.Item(5694, 2).Value = 'Okres Lučenec'
$ws.Cells.Item(5694, 3).Value = 519
$ws.Cells.Item(5694, 4).Value = 6
$ws.Cells.Item(5694, 5).Value = 525
$ws.Cells.Item(5695, 1).Value = 44215
$ws.Cells.Item(5695, 2).Value = 'Okres Malacky'
$ws.Cells.Item(5695, 3).Value = 14
$ws.Cells.Item(5695, 4).Value = 2
$ws.Cells.Item(5695, 5).Value = 16
$ws.Cells.Item(5696, 1).Value = 44215
$ws.Cells.Item(5696, 2).Value = 'Okres Martin'
$ws.Cells.Item(5696, 3).Value = 4383
$ws.Cells.Item(5696, 4).Value = 158
$ws.Cells.Item(5696, 5).Value = 4541
$ws.Cells.Item(5697, 1).Value = 44215
$ws.Cells.Item(5697, 2).Value = 'Okres Michalovce'
$ws.Cells.Item(5697, 3).Value = 460
$ws.Cells.Item(5697, 4).Value = 10
$ws.Cells.Item(5697, 5).Value = 470
$ws.Cells.Item(5698, 1).Value = 44215
$ws.Cells.Item(5698, 2).Value = 'Okres Myjava'
$ws.Cells.Item(5698, 3).Value = 80
$ws.Cells.Item(5698, 4).Value = 0
$ws.Cells.Item(5698, 5).Value = 80
$ws.Cells.Item(5699, 1).Value = 44215
$ws.Cells.Item(5699, 2).Value = 'Okres Námestovo'
$ws.Cells.Item(5699, 3).Value = 14
$ws.Cells.Item(5699, 4).Value = 0
$ws.Cells.Item(5699, 5).Value = 14
$ws.Cells.Item(5700, 1).Value = 44215
$ws.Cells.Item(5700, 2).Value = 'Okres Nitra'
$ws.Cells.Item(5700, 3).Value = 162
$ws.Cells.Item(5700, 4).Value = 9
$ws.Cells.Item(5700, 5).Value = 171
$ws.Cells.Item(5701, 1).Value = 44215
$ws.Cells.Item(5701, 2).Value = 'Okres Nové Mesto nad Váhom'
$ws.Cells.Item(5701, 3).Value = 14
$ws.Cells.Item(5701, 4).Value = 1
$ws.Cells.Item(5701, 5).Value = 15
$ws.Cells.Item(5702, 1).Value = 44215
$ws.Cells.Item(5702, 2).Value = 'Okres Nové Zámky'
$ws.Cells.Item(5702, 3).Value = 3
$ws.Cells.Item(5702, 4).Value = 0
$ws.Cells.Item(5702, 5).Value = 3
$ws.Cells.Item(5703, 1).Value = 44215
$ws.Cells.Item(5703, 2).Value = 'Okres Partizánske'
$ws.Cells.Item(5703, 3).Value = 425
$ws.Cells.Item(5703, 4).Value = 78
$ws.Cells.Item(5703, 5).Value = 503
$ws.Cells.Item(5704, 1).Value = 44215
$ws.Cells.Item(5704, 2).Value = 'Okres Pezinok'
$ws.Cells.Item(5704, 3).Value = 288
$ws.Cells.Item(5704, 4).Value = 10
$ws.Cells.Item(5704, 5).Value = 298
$ws.Cells.Item(5705, 1).Value = 44215
$ws.Cells.Item(5705, 2).Value = 'Okres Piešťany'
$ws.Cells.Item(5705, 3).Value = 700
$ws.Cells.Item(5705, 4).Value = 14
$ws.Cells.Item(5705, 5).Value = 714
$ws.Cells.Item(5706, 1).Value = 44215
$ws.Cells.Item(5706, 2).Value = 'Okres Poprad'
$ws.Cells.Item(5706, 3).Value = 1070
$ws.Cells.Item(5706, 4).Value = 16
$ws.Cells.Item(5706, 5).Value = 1086
$ws.Cells.Item(5707, 1).Value = 44215
$ws.Cells.Item(5707, 2).Value = 'Okres Považská Bystrica'
$ws.Cells.Item(5707, 3).Value = 843
$ws.Cells.Item(5707, 4).Value = 35
$ws.Cells.Item(5707, 5).Value = 878
$ws.Cells.Item(5708, 1).Value = 44215
$ws.Cells.Item(5708, 2).Value = 'Okres Prešov'
$ws.Cells.Item(5708, 3).Value = 342
$ws.Cells.Item(5708, 4).Value = 19
$ws.Cells.Item(5708, 5).Value = 361
$ws.Cells.Item(5709, 1).Value = 44215
$ws.Cells.Item(5709, 2).Value = 'Okres Prievidza'
$ws.Cells.Item(5709, 3).Value = 439
$ws.Cells.Item(5709, 4).Value = 10
$ws.Cells.Item(5709, 5).Value = 449
$ws.Cells.Item(5710, 1).Value = 44215
$ws.Cells.Item(5710, 2).Value = 'Okres Púchov'
$ws.Cells.Item(5710, 3).Value = 0
$ws.Cells.Item(5710, 4).Value = 0
$ws.Cells.Item(5710, 5).Value = 0
$ws.Cells.Item(5711, 1).Value = 44215
$ws.Cells.Item(5711, 2).Value = 'Okres Revúca'
$ws.Cells.Item(5711, 3).Value = 538
$ws.Cells.Item(5711, 4).Value = 49
$ws.Cells.Item(5711, 5).Value = 587
$ws.Cells.Item(5712, 1).Value = 44215
$ws.Cells.Item(5712, 2).Value = 'Okres Rimavská Sobota'
$ws.Cells.Item(5712, 3).Value = 468
$ws.Cells.Item(5712, 4).Value = 21
$ws.Cells.Item(5712, 5).Value = 489
$ws.Cells.Item(5713, 1).Value = 44215
$ws.Cells.Item(5713, 2).Value = 'Okres Rožňava'
$ws.Cells.Item(5713, 3).Value = 555
$ws.Cells.Item(5713, 4).Value = 17
$ws.Cells.Item(5713, 5).Value = 572
$ws.Cells.Item(5714, 1).Value = 44215
$ws.Cells.Item(5714, 2).Value = 'Okres Ružomberok'
$ws.Cells.Item(5714, 3).Value = 1089
$ws.Cells.Item(5714, 4).Value = 7
$ws.Cells.Item(5714, 5).Value = 1096
$ws.Cells.Item(5715, 1).Value = 44215
$ws.Cells.Item(5715, 2).Value = 'Okres Senica'
$ws.Cells.Item(5715, 3).Value = 7
$ws.Cells.Item(5715, 4).Value = 2
$ws.Cells.Item(5715, 5).Value = 9
$ws.Cells.Item(5716, 1).Value = 44215
$ws.Cells.Item(5716, 2).Value = 'Okres Skalica'
$ws.Cells.Item(5716, 3).Value = 459
$ws.Cells.Item(5716, 4).Value = 8
$ws.Cells.Item(5716, 5).Value = 467
$ws.Cells.Item(5717, 1).Value = 44215
$ws.Cells.Item(5717, 2).Value = 'Okres Sobrance'
$ws.Cells.Item(5717, 3).Value = 144
$ws.Cells.Item(5717, 4).Value = 6
$ws.Cells.Item(5717, 5).Value = 150
$ws.Cells.Item(5718, 1).Value = 44215
$ws.Cells.Item(5718, 2).Value = 'Okres Spišská Nová Ves'
$ws.Cells.Item(5718, 3).Value = 145
$ws.Cells.Item(5718, 4).Value = 7
$ws.Cells.Item(5718, 5).Value = 152
$ws.Cells.Item(5719, 1).Value = 44215
$ws.Cells.Item(5719, 2).Value = 'Okres Stará Ľubovňa'
$ws.Cells.Item(5719, 3).Value = 259
$ws.Cells.Item(5719, 4).Value = 7
$ws.Cells.Item(5719, 5).Value = 266
$ws.Cells.Item(5720, 1).Value = 44215
$ws.Cells.Item(5720, 2).Value = 'Okres Stropkov'
$ws.Cells.Item(5720, 3).Value = 2
$ws.Cells.Item(5720, 4).Value = 0
$ws.Cells.Item(5720, 5).Value = 2
$ws.Cells.Item(5721, 1).Value = 44215
$ws.Cells.Item(5721, 2).Value = 'Okres Šaľa'
$ws.Cells.Item(5721, 3).Value = 4
$ws.Cells.Item(5721, 4).Value = 0
$ws.Cells.Item(5721, 5).Value = 4
$ws.Cells.Item(5722, 1).Value = 44215
$ws.Cells.Item(5722, 2).Value = 'Okres Topoľčany'
$ws.Cells.Item(5722, 3).Value = 852
$ws.Cells.Item(5722, 4).Value = 64
$ws.Cells.Item(5722, 5).Value = 916
$ws.Cells.Item(5723, 1).Value = 44215
$ws.Cells.Item(5723, 2).Value = 'Okres Trebišov'
$ws.Cells.Item(5723, 3).Value = 774
$ws.Cells.Item(5723, 4).Value = 21
$ws.Cells.Item(5723, 5).Value = 795
$ws.Cells.Item(5724, 1).Value = 44215
$ws.Cells.Item(5724, 2).Value = 'Okres Trenčín'
$ws.Cells.Item(5724, 3).Value = 571
$ws.Cells.Item(5724, 4).Value = 30
$ws.Cells.Item(5724, 5).Value = 601
$ws.Cells.Item(5725, 1).Value = 44215
$ws.Cells.Item(5725, 2).Value = 'Okres Trnava'
$ws.Cells.Item(5725, 3).Value = 272
$ws.Cells.Item(5725, 4).Value = 12
$ws.Cells.Item(5725, 5).Value = 284
$ws.Cells.Item(5726, 1).Value = 44215
$ws.Cells.Item(5726, 2).Value = 'Okres Turčianske Teplice'
$ws.Cells.Item(5726, 3).Value = 210
$ws.Cells.Item(5726, 4).Value = 2
$ws.Cells.Item(5726, 5).Value = 212
$ws.Cells.Item(5727, 1).Value = 44215
$ws.Cells.Item(5727, 2).Value = 'Okres Tvrdošín'
$ws.Cells.Item(5727, 3).Value = 236
$ws.Cells.Item(5727, 4).Value = 0
$ws.Cells.Item(5727, 5).Value = 236
$ws.Cells.Item(5728, 1).Value = 44215
$ws.Cells.Item(5728, 2).Value = 'Okres Veľký Krtíš'
$ws.Cells.Item(5728, 3).Value = 454
$ws.Cells.Item(5728, 4).Value = 9
$ws.Cells.Item(5728, 5).Value = 463
$ws.Cells.Item(5729, 1).Value = 44215
$ws.Cells.Item(5729, 2).Value = 'Okres Vranov nad Topľou'
$ws.Cells.Item(5729, 3).Value = 335
$ws.Cells.Item(5729, 4).Value = 16
$ws.Cells.Item(5729, 5).Value = 351
$ws.Cells.Item(5730, 1).Value = 44215
$ws.Cells.Item(5730, 2).Value = 'Okres Zvolen'
$ws.Cells.Item(5730, 3).Value = 151
$ws.Cells.Item(5730, 4).Value = 5
$ws.Cells.Item(5730, 5).Value = 156
$ws.Cells.Item(5731, 1).Value = 44215
$ws.Cells.Item(5731, 2).Value = 'Okres Žiar nad Hronom'
$ws.Cells.Item(5731, 3).Value = 355
$ws.Cells.Item(5731, 4).Value = 16
$ws.Cells.Item(5731, 5).Value = 371
$ws.Cells.Item(5732, 1).Value = 44215
$ws.Cells.Item(5732, 2).Value = 'Okres Žilina'
$ws.Cells.Item(5732, 3).Value = 1250
$ws.Cells.Item(5732, 4).Value = 41
$ws.Cells.Item(5732, 5).Value = 1291
